$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the note text for row 22 (shared string "q" -> full note)
$ws.Range("E22").Value = "Returned all materials…"

# Row 21: add Time Out (B21); D21 (duration) recalculates automatically
$ws.Range("B21").Value = 0.72916666666666663

# Row 22: add Time In (A22) and Date (C22); D22 (duration) recalculates automatically
$ws.Range("A22").Value = 0.625
$ws.Range("C22").Value = 42095

# Move the active selection to E22 to match the saved view state
$ws.Range("E22").Select()
